$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Session" to "Anatomy"
$ws.Name = "Anatomy"

# Append 43 new QR-scanner log rows (rows 85-127)
$newRows = @(
    @(85, '244942', 'Anatomy', '05/11/2025', '10:08:41', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(86, '244968', 'Anatomy', '05/11/2025', '10:08:46', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(87, '244938', 'Anatomy', '05/11/2025', '10:08:51', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(88, '244934', 'Anatomy', '05/11/2025', '10:08:53', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(89, '244940', 'Anatomy', '05/11/2025', '10:08:59', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(90, '244925', 'Anatomy', '05/11/2025', '10:09:06', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(91, '234696', 'Anatomy', '05/11/2025', '10:09:32', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(92, '234118', 'Anatomy', '05/11/2025', '10:09:38', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(93, '234128', 'Anatomy', '05/11/2025', '10:09:45', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(94, '244897', 'Anatomy', '05/11/2025', '10:10:16', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(95, '234216', 'Anatomy', '05/11/2025', '10:10:34', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(96, '234094', 'Anatomy', '05/11/2025', '10:10:45', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(97, '234215', 'Anatomy', '05/11/2025', '10:10:57', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(98, '234171', 'Anatomy', '05/11/2025', '10:11:03', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(99, '244956', 'Anatomy', '05/11/2025', '10:11:07', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(100, '244941', 'Anatomy', '05/11/2025', '10:11:08', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(101, '244949', 'Anatomy', '05/11/2025', '10:11:45', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(102, '234347', 'Anatomy', '05/11/2025', '10:11:51', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(103, '234336', 'Anatomy', '05/11/2025', '10:11:56', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(104, '234191', 'Anatomy', '05/11/2025', '10:12:27', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(105, '234193', 'Anatomy', '05/11/2025', '10:12:40', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(106, '234116', 'Anatomy', '05/11/2025', '10:12:46', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(107, '234219', 'Anatomy', '05/11/2025', '10:12:54', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(108, '244902', 'Anatomy', '05/11/2025', '10:13:08', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(109, '234241', 'Anatomy', '05/11/2025', '10:13:24', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(110, '234098', 'Anatomy', '05/11/2025', '10:13:48', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(111, '244896', 'Anatomy', '05/11/2025', '10:20:04', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(112, '244736', 'Anatomy', '05/11/2025', '10:20:11', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(113, '244672', 'Anatomy', '05/11/2025', '10:20:32', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(114, '244841', 'Anatomy', '05/11/2025', '10:22:34', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(115, '244656', 'Anatomy', '05/11/2025', '10:22:59', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(116, '244785', 'Anatomy', '05/11/2025', '10:24:32', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(117, '244756', 'Anatomy', '05/11/2025', '10:24:47', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(118, '244647', 'Anatomy', '05/11/2025', '10:26:58', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(119, '244645', 'Anatomy', '05/11/2025', '10:27:10', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(120, '244893', 'Anatomy', '05/11/2025', '10:28:17', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(121, '244937', 'Anatomy', '05/11/2025', '10:28:24', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(122, '244714', 'Anatomy', '05/11/2025', '10:28:33', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(123, '244765', 'Anatomy', '05/11/2025', '10:32:09', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(124, '244906', 'Anatomy', '05/11/2025', '10:32:46', 'Scan', 'nahla.nagiub@med.asu.edu.eg'),
    @(125, '234829', 'Anatomy', '05/11/2025', '10:38:34', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(126, '234823', 'Anatomy', '05/11/2025', '10:38:57', 'Manual', 'nahla.nagiub@med.asu.edu.eg'),
    @(127, '234830', 'Anatomy', '05/11/2025', '10:39:16', 'Manual', 'nahla.nagiub@med.asu.edu.eg')
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = "'" + $r[1]
    $ws.Cells.Item($rowNum, 2).Value = "'" + $r[2]
    $ws.Cells.Item($rowNum, 3).Value = "'" + $r[3]
    $ws.Cells.Item($rowNum, 4).Value = "'" + $r[4]
    $ws.Cells.Item($rowNum, 5).Value = "'" + $r[5]
    $ws.Cells.Item($rowNum, 6).Value = "'" + $r[6]
}
